# Auto-generated Excel COM-interop script to apply scheduled market-price refresh
# to the Leve profit tables across all job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 1972.375
$ws.Range("I7").Value = 1972.375
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1972.375
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1860.375
$ws.Range("N7").ClearContents()

$ws.Range("H14").Value = 1972.375
$ws.Range("I14").Value = 1972.375
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1972.375
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -1781.375
$ws.Range("N14").ClearContents()

$ws.Range("H41").Value = 248.4375
$ws.Range("I41").Value = 218.41667
$ws.Range("J41").Value = 338.5
$ws.Range("K41").Value = 218.41667
$ws.Range("L41").Value = 338.5
$ws.Range("M41").Value = 221.58333
$ws.Range("N41").Value = -1218.5

$ws.Range("H86").Value = 2790.1936
$ws.Range("I86").Value = 1638.3077
$ws.Range("K86").Value = 1638.3077
$ws.Range("M86").Value = -515.3077000000001

$ws.Range("H89").Value = 2790.1936
$ws.Range("I89").Value = 1638.3077
$ws.Range("K89").Value = 8191.538500000001
$ws.Range("M89").Value = -2575.538500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 15666.667
$ws.Range("J9").Value = 15666.667
$ws.Range("L9").Value = 15666.667
$ws.Range("N9").Value = -16006.667

$ws.Range("H20").Value = 15666.667
$ws.Range("J20").Value = 15666.667
$ws.Range("L20").Value = 15666.667
$ws.Range("N20").Value = -16206.667

$ws.Range("H61").Value = 1801.275
$ws.Range("I61").Value = 1218.6562
$ws.Range("K61").Value = 1218.6562
$ws.Range("M61").Value = -1006.6562

$ws.Range("H88").Value = 4004
$ws.Range("I88").Value = 6500
$ws.Range("J88").Value = 2340
$ws.Range("K88").Value = 6500
$ws.Range("L88").Value = 2340
$ws.Range("M88").Value = -6094
$ws.Range("N88").Value = -3152

$ws.Range("H91").Value = 4004
$ws.Range("I91").Value = 6500
$ws.Range("J91").Value = 2340
$ws.Range("K91").Value = 6500
$ws.Range("L91").Value = 2340
$ws.Range("M91").Value = -5096
$ws.Range("N91").Value = -5148

$ws.Range("H136").Value = 1801.275
$ws.Range("I136").Value = 1218.6562
$ws.Range("K136").Value = 3655.9686
$ws.Range("M136").Value = -1105.9686

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3666.6667
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = -1877
$ws.Range("N86").Value = -6246

$ws.Range("H89").Value = 3666.6667
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = -9384
$ws.Range("N89").Value = -31232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1275613.1
$ws.Range("I6").Value = 1275613.1
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 1275613.1
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -1275500.1
$ws.Range("N6").ClearContents()

$ws.Range("H132").Value = 2484.7693
$ws.Range("I132").Value = 1973.6666
$ws.Range("K132").Value = 5920.9998
$ws.Range("M132").Value = -3390.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 1980
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H122").Value = 10753825
$ws.Range("I122").Value = 18518884
$ws.Range("J122").Value = 2203.7693
$ws.Range("K122").Value = 166669956
$ws.Range("L122").Value = 19833.9237
$ws.Range("M122").Value = -166667506
$ws.Range("N122").Value = -24733.9237

$ws.Range("H132").Value = 37037948
$ws.Range("I132").Value = 55556444
$ws.Range("J132").Value = 956.3333
$ws.Range("K132").Value = 500007996
$ws.Range("L132").Value = 8606.9997
$ws.Range("M132").Value = -500005466
$ws.Range("N132").Value = -13666.9997

$ws.Range("H133").Value = 4363.364
$ws.Range("I133").Value = 4204
$ws.Range("J133").Value = 4391.8213
$ws.Range("K133").Value = 12612
$ws.Range("L133").Value = 13175.4639
$ws.Range("M133").Value = -7552
$ws.Range("N133").Value = -23295.4639

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1806.2354
$ws.Range("I113").Value = 1853
$ws.Range("J113").Value = 1720.5
$ws.Range("K113").Value = 1853
$ws.Range("L113").Value = 1720.5
$ws.Range("M113").Value = 317
$ws.Range("N113").Value = -6060.5

$ws.Range("H122").Value = 5755.727
$ws.Range("I122").Value = 5538.625
$ws.Range("J122").Value = 6334.6665
$ws.Range("K122").Value = 16615.875
$ws.Range("L122").Value = 19003.9995
$ws.Range("M122").Value = -14165.875
$ws.Range("N122").Value = -23903.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3665.3333
$ws.Range("I16").Value = 4073.5
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 4073.5
$ws.Range("L16").Value = 400
$ws.Range("M16").Value = -3903.5
$ws.Range("N16").Value = -740

$ws.Range("H82").Value = 1398.15
$ws.Range("I82").Value = 1174.6923
$ws.Range("K82").Value = 1174.6923
$ws.Range("M82").Value = -813.6922999999999

$ws.Range("H85").Value = 1398.15
$ws.Range("I85").Value = 1174.6923
$ws.Range("K85").Value = 1174.6923
$ws.Range("M85").Value = 73.30770000000007

$ws.Range("H122").Value = 9455.444
$ws.Range("I122").Value = 12200.728
$ws.Range("J122").Value = 5141.4287
$ws.Range("K122").Value = 36602.18399999999
$ws.Range("L122").Value = 15424.2861
$ws.Range("M122").Value = -34152.18399999999
$ws.Range("N122").Value = -20324.2861

$ws.Range("H136").Value = 2522.625
$ws.Range("I136").Value = 1608.4706
$ws.Range("K136").Value = 4825.4118
$ws.Range("M136").Value = -2275.4118

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1209.6666
$ws.Range("I122").Value = 1133.8667
$ws.Range("J122").Value = 1399.1666
$ws.Range("K122").Value = 3401.6001
$ws.Range("L122").Value = 4197.4998
$ws.Range("M122").Value = -951.6001000000001
$ws.Range("N122").Value = -9097.4998

$ws.Range("H126").Value = 3227.35
$ws.Range("I126").Value = 2619
$ws.Range("J126").Value = 5052.4
$ws.Range("K126").Value = 7857
$ws.Range("L126").Value = 15157.2
$ws.Range("M126").Value = -5387
$ws.Range("N126").Value = -20097.2

$ws.Range("H136").Value = 16754.385
$ws.Range("I136").Value = 21480.7
$ws.Range("K136").Value = 64442.10000000001
$ws.Range("M136").Value = -61892.10000000001

